$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns X1 and Y1 were mistakenly labeled as "current year" headers
# (duplicating V1/W1). Relabel them as "prior year" headers instead.
$ws.Range("X1").Value = "Prior Year Dwelling MV"
$ws.Range("Y1").Value = "Prior Year Dwelling Total"
